$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "IF first dice roll of the game, display: Roll some dice!",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "IF first dice roll of the turn, display: Roll some dice!", 2
)
